$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (header row 1, data rows below).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()
$lastCol = $usedRange.Columns.Count()

$firstDataRow = 2

# Read all data rows (time, AA1, AA2, AA3) into an array of rows.
$rows = @()
for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value()
    }
    $rows += , $rowVals
}

# Sort the rows by the first column (time, ascending) - this is the
# calibration timestamp ordering produced after re-running the needle
# calibration.
$sorted = $rows | Sort-Object { $_[0] }

# Write the sorted rows back into the sheet.
for ($i = 0; $i -lt $sorted.Count; $i++) {
    $r = $firstDataRow + $i
    $rowVals = $sorted[$i]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
